$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 116, shifting existing rows 116-122 down to 117-123
$ws.Rows.Item(116).Insert()

# Populate the new row 116 with the new weekly data point
$ws.Cells.Item(116, 1).Value = 6
$ws.Cells.Item(116, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(116, 3).Value = "Metropolitana"
$ws.Cells.Item(116, 4).Value = 44516
$ws.Cells.Item(116, 5).Value = 13
$ws.Cells.Item(116, 6).Value = 100112029
$ws.Cells.Item(116, 7).Value = "Orégano"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 34
$ws.Cells.Item(116, 11).Value = 8000
$ws.Cells.Item(116, 12).Value = 9000
$ws.Cells.Item(116, 13).Value = 8441
$ws.Cells.Item(116, 14).Value = "`$/docena de atados"
$ws.Cells.Item(116, 15).Value = "Región Metropolitana"
$ws.Cells.Item(116, 16).Value = 2814
$ws.Cells.Item(116, 17).Value = 3
$ws.Cells.Item(116, 18).Value = "Hortaliza"
